# Reorder the "Recorded By" (column G) values on the "Session Analysis Results"
# sheet so that any "system"-like tokens (case-insensitive match to "system")
# appear first (preserving their relative order), followed by the remaining
# tokens (e.g. email addresses), also preserving their relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Column G is the "Recorded By" column (G1 header = "Recorded By")
$col = 7

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ",\s*" | ForEach-Object { $_.Trim() }

        $systemParts = @()
        $otherParts = @()
        foreach ($part in $parts) {
            if ($part.ToLower() -eq "system") {
                $systemParts += $part
            } else {
                $otherParts += $part
            }
        }

        if ($systemParts.Count -gt 0) {
            $newParts = $systemParts + $otherParts
            $newValue = [string]::Join(", ", $newParts)
            if ($newValue -ne $value) {
                $cell.Value2 = $newValue
            }
        }
    }
}
